$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = '40236112'
$ws.Cells.Item(2, 2).Value = 'https://oleks-netizen.github.io/product-images/40236112/1.jpg,https://oleks-netizen.github.io/product-images/40236112/7.jpg,https://oleks-netizen.github.io/product-images/40236112/2.jpg,https://oleks-netizen.github.io/product-images/40236112/4.jpg,https://oleks-netizen.github.io/product-images/40236112/5.jpg,https://oleks-netizen.github.io/product-images/40236112/6.jpg'
$ws.Cells.Item(2, 3).Value = 6

$ws.Cells.Item(3, 1).Value = 'BN-GC-14-1-o-felt-d-Photoroom'
$ws.Cells.Item(3, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-GC-14-1-o-felt-d-Photoroom/1.jpg,https://oleks-netizen.github.io/product-images/BN-GC-14-1-o-felt-d-Photoroom/2.jpg,https://oleks-netizen.github.io/product-images/BN-GC-14-1-o-felt-d-Photoroom/3.jpg'
$ws.Cells.Item(3, 3).Value = 3

$ws.Cells.Item(4, 1).Value = 'BN-GC-14-1-o-felt-d'
$ws.Cells.Item(4, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-GC-14-1-o-felt-d/6.jpg'
$ws.Cells.Item(4, 3).Value = 1

$ws.Cells.Item(5, 1).Value = 'BN-GC-16-1-o-felt-d-Photoroom'
$ws.Cells.Item(5, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-GC-16-1-o-felt-d-Photoroom/1.jpg,https://oleks-netizen.github.io/product-images/BN-GC-16-1-o-felt-d-Photoroom/2.jpg,https://oleks-netizen.github.io/product-images/BN-GC-16-1-o-felt-d-Photoroom/3.jpg'
$ws.Cells.Item(5, 3).Value = 3

$ws.Cells.Item(6, 1).Value = 'BN-GC-16-1-o-felt-d'
$ws.Cells.Item(6, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-GC-16-1-o-felt-d/6.jpg'
$ws.Cells.Item(6, 3).Value = 1

$ws.Cells.Item(7, 1).Value = 'BN-OP-12-g-kr'
$ws.Cells.Item(7, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-OP-12-g-kr/1.jpg'
$ws.Cells.Item(7, 3).Value = 1

$ws.Cells.Item(8, 1).Value = 'BN-OP-12-g-kr-Photoroom'
$ws.Cells.Item(8, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-OP-12-g-kr-Photoroom/2.jpg,https://oleks-netizen.github.io/product-images/BN-OP-12-g-kr-Photoroom/3.jpg'
$ws.Cells.Item(8, 3).Value = 2

$ws.Cells.Item(9, 1).Value = 'BN-SB-13-k'
$ws.Cells.Item(9, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-SB-13-k/1.jpg'
$ws.Cells.Item(9, 3).Value = 1

$ws.Cells.Item(10, 1).Value = 'BN-SB-13-k-Photoroom'
$ws.Cells.Item(10, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-SB-13-k-Photoroom/2.jpg,https://oleks-netizen.github.io/product-images/BN-SB-13-k-Photoroom/3.jpg'
$ws.Cells.Item(10, 3).Value = 2

$ws.Cells.Item(11, 1).Value = 'BN-SB-13-vin'
$ws.Cells.Item(11, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-SB-13-vin/1.jpg'
$ws.Cells.Item(11, 3).Value = 1

$ws.Cells.Item(12, 1).Value = 'BN-SB-13-vin-Photoroom'
$ws.Cells.Item(12, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-SB-13-vin-Photoroom/2.jpg,https://oleks-netizen.github.io/product-images/BN-SB-13-vin-Photoroom/3.jpg'
$ws.Cells.Item(12, 3).Value = 2

$ws.Cells.Item(13, 1).Value = 'BN-SB-6_073926907251'
$ws.Cells.Item(13, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-SB-6_073926907251/1.jpg'
$ws.Cells.Item(13, 3).Value = 1

$ws.Cells.Item(14, 1).Value = 'BN-SB-6_073926337432-Photoroom'
$ws.Cells.Item(14, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-SB-6_073926337432-Photoroom/2.jpg'
$ws.Cells.Item(14, 3).Value = 1

$ws.Cells.Item(15, 1).Value = 'BN-SB-6_073926493725-Photoroom'
$ws.Cells.Item(15, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-SB-6_073926493725-Photoroom/3.jpg'
$ws.Cells.Item(15, 3).Value = 1

$ws.Cells.Item(16, 1).Value = 'BN-SB-6_073926689097-Photoroom'
$ws.Cells.Item(16, 2).Value = 'https://oleks-netizen.github.io/product-images/BN-SB-6_073926689097-Photoroom/4.jpg'
$ws.Cells.Item(16, 3).Value = 1

$ws.Cells.Item(17, 1).Value = 'HB10032_3044'
$ws.Cells.Item(17, 2).Value = 'https://oleks-netizen.github.io/product-images/HB10032_3044/1.jpg,https://oleks-netizen.github.io/product-images/HB10032_3044/5.jpg,https://oleks-netizen.github.io/product-images/HB10032_3044/3.jpg,https://oleks-netizen.github.io/product-images/HB10032_3044/4.jpg'
$ws.Cells.Item(17, 3).Value = 4

$ws.Cells.Item(18, 1).Value = 'RB63 LIME M'
$ws.Cells.Item(18, 2).Value = 'https://oleks-netizen.github.io/product-images/RB63 LIME M/1.jpg,https://oleks-netizen.github.io/product-images/RB63 LIME M/5.jpg,https://oleks-netizen.github.io/product-images/RB63 LIME M/2.jpg,https://oleks-netizen.github.io/product-images/RB63 LIME M/3.jpg'
$ws.Cells.Item(18, 3).Value = 4

$ws.Cells.Item(19, 1).Value = 'TW-PH-beige-ksr'
$ws.Cells.Item(19, 2).Value = 'https://oleks-netizen.github.io/product-images/TW-PH-beige-ksr/1.jpg,https://oleks-netizen.github.io/product-images/TW-PH-beige-ksr/3.jpg'
$ws.Cells.Item(19, 3).Value = 2

$ws.Cells.Item(20, 1).Value = 'TW-PH-black-ksr-Photoroom'
$ws.Cells.Item(20, 2).Value = 'https://oleks-netizen.github.io/product-images/TW-PH-black-ksr-Photoroom/1.jpg,https://oleks-netizen.github.io/product-images/TW-PH-black-ksr-Photoroom/4.jpg'
$ws.Cells.Item(20, 3).Value = 2

$ws.Cells.Item(21, 1).Value = 'TW-PH-black-ksr'
$ws.Cells.Item(21, 2).Value = 'https://oleks-netizen.github.io/product-images/TW-PH-black-ksr/6.jpg,https://oleks-netizen.github.io/product-images/TW-PH-black-ksr/3.jpg,https://oleks-netizen.github.io/product-images/TW-PH-black-ksr/5.jpg'
$ws.Cells.Item(21, 3).Value = 3

$ws.Cells.Item(22, 1).Value = 'TW-PH-dark-blue'
$ws.Cells.Item(22, 2).Value = 'https://oleks-netizen.github.io/product-images/TW-PH-dark-blue/1.jpg,https://oleks-netizen.github.io/product-images/TW-PH-dark-blue/3.jpg'
$ws.Cells.Item(22, 3).Value = 2

$ws.Cells.Item(23, 1).Value = 'TW-PH-kon-crz'
$ws.Cells.Item(23, 2).Value = 'https://oleks-netizen.github.io/product-images/TW-PH-kon-crz/1.jpg,https://oleks-netizen.github.io/product-images/TW-PH-kon-crz/3.jpg'
$ws.Cells.Item(23, 3).Value = 2

$ws.Cells.Item(24, 1).Value = 'TW-PH-kon-ksr'
$ws.Cells.Item(24, 2).Value = 'https://oleks-netizen.github.io/product-images/TW-PH-kon-ksr/1.jpg,https://oleks-netizen.github.io/product-images/TW-PH-kon-ksr/3.jpg'
$ws.Cells.Item(24, 3).Value = 2

$ws.Cells.Item(25, 1).Value = 'TW-PH-mars-ksr'
$ws.Cells.Item(25, 2).Value = 'https://oleks-netizen.github.io/product-images/TW-PH-mars-ksr/1.jpg,https://oleks-netizen.github.io/product-images/TW-PH-mars-ksr/3.jpg'
$ws.Cells.Item(25, 3).Value = 2

$ws.Cells.Item(26, 1).Value = 'TW-PH-red-saf'
$ws.Cells.Item(26, 2).Value = 'https://oleks-netizen.github.io/product-images/TW-PH-red-saf/1.jpg,https://oleks-netizen.github.io/product-images/TW-PH-red-saf/3.jpg'
$ws.Cells.Item(26, 3).Value = 2

